$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 194; this shifts the existing rows 194-204
# (and all their data) down to 195-205, matching the rest of the diff
# automatically.
$ws.Rows(194).Insert()

# Populate the newly inserted row 194 with a new weekly data point for
# "Feria Lagunitas de Puerto Montt" / Apio / Americana (o) / Primera.
$ws.Range("A194").Value = 4
$ws.Range("B194").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C194").Value = "Los Lagos"
$ws.Range("D194").Value = 44610
$ws.Range("E194").Value = 10
$ws.Range("F194").Value = 100112017
$ws.Range("G194").Value = "Apio"
$ws.Range("H194").Value = "Americana (o)"
$ws.Range("I194").Value = "Primera"
$ws.Range("J194").Value = 40
$ws.Range("K194").Value = 12000
$ws.Range("L194").Value = 12000
$ws.Range("M194").Value = 12000
$ws.Range("N194").Value = "$/docena de matas"
$ws.Range("O194").Value = "Región de Coquimbo"
$ws.Range("P194").Value = 2000
$ws.Range("Q194").Value = 6
$ws.Range("R194").Value = "Hortaliza"
